$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.344.69'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.902.04'
$ws.Range('E3').Value = '  +2.20%  '
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.77%  '
$ws.Range('E6').Value = '  +6.10%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.40'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.18%  '
$ws.Range('E9').Value = '  +5.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.99'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0721'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0994'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').Value = '2.179.73'
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.697'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.99%  '
$ws.Range('D16').Value = '1.906.65'
$ws.Range('E16').Value = '  +2.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.56%  '
$ws.Range('D18').Value = '35.326.00'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.15%  '
$ws.Range('D20').Value = '0.0₃0827'
$ws.Range('E20').Value = '  +4.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '240.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.49'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +24.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.41'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.36'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.52%  '
$ws.Range('E30').Value = '  +2.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0564'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.932'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.09'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.56%  '
$ws.Range('E36').Value = '  -4.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.03'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.29%  '
$ws.Range('E38').Value = '  +1.79%  '
$ws.Range('E39').Value = '  -0.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0207'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.29%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '89.56'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0622'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.18%  '
$ws.Range('D44').Value = '1.338.04'
$ws.Range('E44').Value = '  -0.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.39'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '48.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +39.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.78'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.69%  '
$ws.Range('D51').Value = '2.089.79'
$ws.Range('E51').Value = '  +2.23%  '
